# Weekly update: insert a new daily price record for Frambuesa
# (Vega Central Mapocho de Santiago) ahead of the existing rows.
# This shifts all rows from 152..169 down to 153..170 and fills the
# freshly inserted row 152 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 152, pushing existing data down.
$ws.Rows.Item(152).Insert()

# Populate the newly inserted row with the new record.
$ws.Cells.Item(152, 1).Value = 9
$ws.Cells.Item(152, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(152, 3).Value = "Metropolitana"
$ws.Cells.Item(152, 4).Value = 45077
$ws.Cells.Item(152, 5).Value = 13
$ws.Cells.Item(152, 6).Value = "Fruta"
$ws.Cells.Item(152, 7).Value = 100101
$ws.Cells.Item(152, 8).Value = "Berries"
$ws.Cells.Item(152, 9).Value = 100101004
$ws.Cells.Item(152, 10).Value = "Frambuesa"
$ws.Cells.Item(152, 11).Value = "Sin especificar"
$ws.Cells.Item(152, 12).Value = "Primera"
$ws.Cells.Item(152, 13).Value = 500
$ws.Cells.Item(152, 14).Value = 9000
$ws.Cells.Item(152, 15).Value = 9500
$ws.Cells.Item(152, 16).Value = 9280
$ws.Cells.Item(152, 17).Value = "`$/bandeja 2 kilos"
$ws.Cells.Item(152, 18).Value = "Provincia de Linares"
$ws.Cells.Item(152, 19).Value = 4640
$ws.Cells.Item(152, 20).Value = 2
